$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 8302.735296243894
$ws.Range("C2").Value = 1750.482923210865
$ws.Range("D2").Value = 3186.185724672492
